$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free approach: set text-looking numeric/percent values explicitly as text
# by forcing the cell NumberFormat to Text ("@") before assignment, then resetting
# the style back to "Normal" so no stray number-format style lingers on the cell.

# Row 2
$r = $ws.Range("D2")
$r.NumberFormat = "@"
$r.Value = '334.94'
$r.Style = "Normal"

# Row 3
$r = $ws.Range("D3")
$r.NumberFormat = "@"
$r.Value = '43.86'
$r.Style = "Normal"
$r = $ws.Range("E3")
$r.NumberFormat = "@"
$r.Value = '6.24%'
$r.Style = "Normal"

# Row 4
$r = $ws.Range("D4")
$r.NumberFormat = "@"
$r.Value = '5.727'
$r.Style = "Normal"
$r = $ws.Range("E4")
$r.NumberFormat = "@"
$r.Value = '1.80%'
$r.Style = "Normal"

# Row 5
$r = $ws.Range("D5")
$r.NumberFormat = "@"
$r.Value = '0.08330'
$r.Style = "Normal"
$r = $ws.Range("E5")
$r.NumberFormat = "@"
$r.Value = '1.38%'
$r.Style = "Normal"

# Row 6
$r = $ws.Range("D6")
$r.NumberFormat = "@"
$r.Value = '8.850'
$r.Style = "Normal"
$r = $ws.Range("E6")
$r.NumberFormat = "@"
$r.Value = '1.09%'
$r.Style = "Normal"

# Row 7
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$r = $ws.Range("D7")
$r.NumberFormat = "@"
$r.Value = '4.516'
$r.Style = "Normal"
$r = $ws.Range("E7")
$r.NumberFormat = "@"
$r.Value = '0.28%'
$r.Style = "Normal"

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$r = $ws.Range("D8")
$r.NumberFormat = "@"
$r.Value = '1.960'
$r.Style = "Normal"
$r = $ws.Range("E8")
$r.NumberFormat = "@"
$r.Value = '-2.48%'
$r.Style = "Normal"

# Row 9
$ws.Range("B9").Value = 'BTSEToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$r = $ws.Range("D9")
$r.NumberFormat = "@"
$r.Value = '2.871'
$r.Style = "Normal"
$r = $ws.Range("E9")
$r.NumberFormat = "@"
$r.Value = '-3.87%'
$r.Style = "Normal"

# Row 10
$ws.Range("B10").Value = 'MXToken'
$ws.Range("C10").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$r = $ws.Range("D10")
$r.NumberFormat = "@"
$r.Value = '0.9439'
$r.Style = "Normal"
$r = $ws.Range("E10")
$r.NumberFormat = "@"
$r.Value = '2.31%'
$r.Style = "Normal"

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$r = $ws.Range("D11")
$r.NumberFormat = "@"
$r.Value = '0.1247'
$r.Style = "Normal"
$r = $ws.Range("E11")
$r.NumberFormat = "@"
$r.Value = '-2.40%'
$r.Style = "Normal"

# Row 12
$ws.Range("B12").Value = 'WazirX'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$r = $ws.Range("D12")
$r.NumberFormat = "@"
$r.Value = '0.1983'
$r.Style = "Normal"
$r = $ws.Range("E12")
$r.NumberFormat = "@"
$r.Value = '1.20%'
$r.Style = "Normal"

# Row 13
$ws.Range("B13").Value = 'MandalaExchangeToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$r = $ws.Range("D13")
$r.NumberFormat = "@"
$r.Value = '0.1070'
$r.Style = "Normal"
$r = $ws.Range("E13")
$r.NumberFormat = "@"
$r.Value = '14.07%'
$r.Style = "Normal"

# Row 14
$ws.Range("B14").Value = 'BitrueCoin'
$ws.Range("C14").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$r = $ws.Range("D14")
$r.NumberFormat = "@"
$r.Value = '0.04538'
$r.Style = "Normal"
$r = $ws.Range("E14")
$r.NumberFormat = "@"
$r.Value = '17.98%'
$r.Style = "Normal"

# Row 15
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$r = $ws.Range("D15")
$r.NumberFormat = "@"
$r.Value = '0.1066'
$r.Style = "Normal"
$r = $ws.Range("E15")
$r.NumberFormat = "@"
$r.Value = '0.67%'
$r.Style = "Normal"

# Row 16
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$r = $ws.Range("D16")
$r.NumberFormat = "@"
$r.Value = '0.001293'
$r.Style = "Normal"
$r = $ws.Range("E16")
$r.NumberFormat = "@"
$r.Value = '-0.95%'
$r.Style = "Normal"

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$r = $ws.Range("D17")
$r.NumberFormat = "@"
$r.Value = '0.005911'
$r.Style = "Normal"
$r = $ws.Range("E17")
$r.NumberFormat = "@"
$r.Value = '-5.10%'
$r.Style = "Normal"

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$r = $ws.Range("D18")
$r.NumberFormat = "@"
$r.Value = '3.497'
$r.Style = "Normal"
$r = $ws.Range("E18")
$r.NumberFormat = "@"
$r.Value = '1.39%'
$r.Style = "Normal"

# Row 20
$r = $ws.Range("D20")
$r.NumberFormat = "@"
$r.Value = '8.663'
$r.Style = "Normal"
$r = $ws.Range("E20")
$r.NumberFormat = "@"
$r.Value = '4.56%'
$r.Style = "Normal"

# Row 21
$r = $ws.Range("D21")
$r.NumberFormat = "@"
$r.Value = '0.1353'
$r.Style = "Normal"
$r = $ws.Range("E21")
$r.NumberFormat = "@"
$r.Value = '-0.94%'
$r.Style = "Normal"

# Row 22
$r = $ws.Range("D22")
$r.NumberFormat = "@"
$r.Value = '0.2692'
$r.Style = "Normal"
$r = $ws.Range("E22")
$r.NumberFormat = "@"
$r.Value = '1.13%'
$r.Style = "Normal"

# Row 23
$r = $ws.Range("D23")
$r.NumberFormat = "@"
$r.Value = '0.04423'
$r.Style = "Normal"
$r = $ws.Range("E23")
$r.NumberFormat = "@"
$r.Value = '0.26%'
$r.Style = "Normal"

# Row 24
$r = $ws.Range("D24")
$r.NumberFormat = "@"
$r.Value = '0.001256'
$r.Style = "Normal"
$r = $ws.Range("E24")
$r.NumberFormat = "@"
$r.Value = '-0.21%'
$r.Style = "Normal"

# Row 25
$r = $ws.Range("D25")
$r.NumberFormat = "@"
$r.Value = '0.004345'
$r.Style = "Normal"
$r = $ws.Range("E25")
$r.NumberFormat = "@"
$r.Value = '0.73%'
$r.Style = "Normal"

# Row 26
$r = $ws.Range("D26")
$r.NumberFormat = "@"
$r.Value = '0.0001261'
$r.Style = "Normal"
$r = $ws.Range("E26")
$r.NumberFormat = "@"
$r.Value = '5.05%'
$r.Style = "Normal"

# Row 39
$r = $ws.Range("D39")
$r.NumberFormat = "@"
$r.Value = '0.02811'
$r.Style = "Normal"
$r = $ws.Range("E39")
$r.NumberFormat = "@"
$r.Value = '2.33%'
$r.Style = "Normal"

# Row 40
$r = $ws.Range("D40")
$r.NumberFormat = "@"
$r.Value = '0.06007'
$r.Style = "Normal"
$r = $ws.Range("E40")
$r.NumberFormat = "@"
$r.Value = '10.04%'
$r.Style = "Normal"

# Row 41
$r = $ws.Range("D41")
$r.NumberFormat = "@"
$r.Value = '0.007931'
$r.Style = "Normal"
$r = $ws.Range("E41")
$r.NumberFormat = "@"
$r.Value = '1.78%'
$r.Style = "Normal"

# Row 42
$r = $ws.Range("D42")
$r.NumberFormat = "@"
$r.Value = '0.1427'
$r.Style = "Normal"
$r = $ws.Range("E42")
$r.NumberFormat = "@"
$r.Value = '0.39%'
$r.Style = "Normal"

# Row 43
$r = $ws.Range("E43")
$r.NumberFormat = "@"
$r.Value = '0.30%'
$r.Style = "Normal"

# Row 44
$r = $ws.Range("D44")
$r.NumberFormat = "@"
$r.Value = '0.002171'
$r.Style = "Normal"
$r = $ws.Range("E44")
$r.NumberFormat = "@"
$r.Value = '0.00%'
$r.Style = "Normal"

# Row 45
$r = $ws.Range("D45")
$r.NumberFormat = "@"
$r.Value = '0.01013'
$r.Style = "Normal"
$r = $ws.Range("E45")
$r.NumberFormat = "@"
$r.Value = '-11.47%'
$r.Style = "Normal"

# Row 46
$r = $ws.Range("E46")
$r.NumberFormat = "@"
$r.Value = '3.44%'
$r.Style = "Normal"

# Row 47
$r = $ws.Range("E47")
$r.NumberFormat = "@"
$r.Value = '0.04%'
$r.Style = "Normal"

# Row 48
$r = $ws.Range("D48")
$r.NumberFormat = "@"
$r.Value = '0.003193'
$r.Style = "Normal"
$r = $ws.Range("E48")
$r.NumberFormat = "@"
$r.Value = '0.06%'
$r.Style = "Normal"

# Row 49
$r = $ws.Range("E49")
$r.NumberFormat = "@"
$r.Value = '-0.34%'
$r.Style = "Normal"

# Row 50
$r = $ws.Range("E50")
$r.NumberFormat = "@"
$r.Value = '0.04%'
$r.Style = "Normal"

# Row 51
$r = $ws.Range("D51")
$r.NumberFormat = "@"
$r.Value = '0.0002002'
$r.Style = "Normal"
$r = $ws.Range("E51")
$r.NumberFormat = "@"
$r.Value = '0.04%'
$r.Style = "Normal"
